$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that are removed ---
$ws.Range("E7").ClearContents()

# --- Set new / changed cell values ---
$ws.Range("A10").Value = "combat simulation"
$ws.Range("B10").Value = "combat simulator"
$ws.Range("A11").Value = "flight simulation"
$ws.Range("B11").Value = "flight simulator"
$ws.Range("A26").Value = "space simulation"
$ws.Range("B26").Value = "space simulator"
$ws.Range("A32").Value = "mecha simulation"
$ws.Range("B32").Value = "mech simulator"
$ws.Range("A41").Value = "real-time strategy"
$ws.Range("B41").Value = "real time strategy"
$ws.Range("A42").Value = "weapon-based fighting"
$ws.Range("B42").Value = "weapon based fighting"
$ws.Range("C42").Value = "2d weapon based fighting"
$ws.Range("A43").Value = "puzzle card"
$ws.Range("B43").Value = "puzzle & cards"
$ws.Range("A44").Value = "flight combat"
$ws.Range("B44").Value = "aerial combat"
$ws.Range("A45").Value = "war"
$ws.Range("B45").Value = "wargame"
$ws.Range("A46").Value = "multiplayer battle arena"
$ws.Range("B46").Value = "moba"
$ws.Range("A47").Value = "fvm"
$ws.Range("B47").Value = "full motion video (fmv)"
$ws.Range("A48").Value = "sci-fi puzzle platform game"
$ws.Range("B48").Value = "sci-fi puzzle-platform game"
$ws.Range("A49").Value = "american-football"
$ws.Range("B49").Value = "american football"
$ws.Range("A50").Value = "soccer"
$ws.Range("B50").Value = "traditional soccer"
$ws.Range("A51").Value = "space combat simulation"
$ws.Range("B51").Value = "space combat sim"
$ws.Range("A52").Value = "3d"
$ws.Range("B52").Value = "3-d"
$ws.Range("A53").Value = "virtual reality"
$ws.Range("B53").Value = "augmented reality"
$ws.Range("A54").Value = "light-gun"
$ws.Range("B54").Value = "light gun"
$ws.Range("A55").Value = "historical"
$ws.Range("B55").Value = "history"
$ws.Range("A56").Value = "historical action adventure"
$ws.Range("B56").Value = "historic action-adventure"
$ws.Range("A57").Value = "third-person action adventure"
$ws.Range("B57").Value = "third-person action-adventure"
$ws.Range("A58").Value = "run gun"
$ws.Range("B58").Value = "run and gun"
$ws.Range("A59").Value = "construction management simulation"
$ws.Range("B59").Value = "construction and management simulation"
$ws.Range("A60").Value = "trivia"
$ws.Range("B60").Value = "quiz"

# --- Apply style (format) to new cells: copy format from existing template cells ---
$ws.Range("B39").Copy()
$ws.Range("B42,C42,B43,B44,B45,B46,B47,B48,B49,B50,B51,B54,B55,A56,B56,A57,B57,B58,A59,B59,A60,B60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update sheet view (scroll position / selection) ---
$ws.Range("B60").Select()
$excel.ActiveWindow.ScrollRow = 34
